# Auto-generated: apply cached-value updates (scheduled market-price refresh)
# to the Leve profit tracker sheets. Mirrors the upstream diff exactly:
# set/clear individual <c> numeric cells per row without touching formulas,
# styles, or any other sheet content.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 4990.1665
$ws.Range("I34").Value = 4990.1665
$ws.Range("K34").Value = 4990.1665
$ws.Range("M34").Value = -4787.1665
$ws.Range("H36").Value = 4990.1665
$ws.Range("I36").Value = 4990.1665
$ws.Range("K36").Value = 4990.1665
$ws.Range("M36").Value = -4275.1665
$ws.Range("H92").Value = 102616.2
$ws.Range("J92").Value = 2995
$ws.Range("L92").Value = 2995
$ws.Range("N92").Value = -5491
$ws.Range("H107").Value = 1439.6
$ws.Range("I107").Value = 1096
$ws.Range("K107").Value = 1096
$ws.Range("M107").Value = 824
$ws.Range("H112").Value = 1487
$ws.Range("J112").Value = 1550.7
$ws.Range("L112").Value = 4652.1
$ws.Range("N112").Value = -6868.1
$ws.Range("H132").Value = 6746.346
$ws.Range("I132").Value = 7359.909
$ws.Range("K132").Value = 22079.727
$ws.Range("M132").Value = -19549.727
$ws.Range("H138").Value = 2728.157
$ws.Range("I138").Value = 2297.8
$ws.Range("K138").Value = 6893.400000000001
$ws.Range("M138").Value = -1753.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9788
$ws.Range("H97").Value = 1132.5
$ws.Range("J97").Value = 1323.5
$ws.Range("L97").Value = 1323.5
$ws.Range("N97").Value = -2315.5
$ws.Range("H102").Value = 1892
$ws.Range("J102").Value = 2159
$ws.Range("L102").Value = 2159
$ws.Range("N102").Value = -5403
$ws.Range("H132").Value = 205379
$ws.Range("I132").Value = 205379
$ws.Range("K132").Value = 616137
$ws.Range("M132").Value = -613607
$ws.Range("H136").Value = 10000
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -27450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4491.6665
$ws.Range("I94").Value = 4487.5
$ws.Range("K94").Value = 4487.5
$ws.Range("M94").Value = -4036.5
$ws.Range("H99").Value = 4103.4185
$ws.Range("I99").Value = 3662.425
$ws.Range("J99").Value = 9983.333000000001
$ws.Range("K99").Value = 3662.425
$ws.Range("L99").Value = 9983.333000000001
$ws.Range("M99").Value = -2164.425
$ws.Range("N99").Value = -12979.333
$ws.Range("H134").Value = 3935.1538
$ws.Range("I134").Value = 3378
$ws.Range("J134").Value = 6999.5
$ws.Range("K134").Value = 10134
$ws.Range("L134").Value = 20998.5
$ws.Range("M134").Value = -7599
$ws.Range("N134").Value = -26068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1907.3846
$ws.Range("I16").Value = 1899.75
$ws.Range("K16").Value = 1899.75
$ws.Range("M16").Value = -1612.75
$ws.Range("H58").Value = 114198.664
$ws.Range("I58").Value = 202159.8
$ws.Range("J58").Value = 4247.25
$ws.Range("K58").Value = 202159.8
$ws.Range("L58").Value = 4247.25
$ws.Range("M58").Value = -201956.8
$ws.Range("N58").Value = -4653.25
$ws.Range("H74").Value = 47584.855
$ws.Range("J74").Value = 47584.855
$ws.Range("L74").Value = 47584.855
$ws.Range("N74").Value = -49332.855
$ws.Range("H77").Value = 47584.855
$ws.Range("J77").Value = 47584.855
$ws.Range("L77").Value = 142754.565
$ws.Range("N77").Value = -151490.565
$ws.Range("H113").Value = 1907.3846
$ws.Range("I113").Value = 1899.75
$ws.Range("K113").Value = 1899.75
$ws.Range("M113").Value = 270.25
$ws.Range("H122").Value = 2172.25
$ws.Range("J122").Value = 2776
$ws.Range("L122").Value = 8328
$ws.Range("N122").Value = -13228
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 114198.664
$ws.Range("I136").Value = 202159.8
$ws.Range("J136").Value = 4247.25
$ws.Range("K136").Value = 606479.3999999999
$ws.Range("L136").Value = 12741.75
$ws.Range("M136").Value = -603929.3999999999
$ws.Range("N136").Value = -17841.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 490
$ws.Range("I5").Value = 490
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1470
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1358
$ws.Range("N5").ClearContents()
$ws.Range("H18").Value = 998
$ws.Range("I18").Value = 998
$ws.Range("K18").Value = 2994
$ws.Range("M18").Value = -2825
$ws.Range("H63").Value = 2250
$ws.Range("I63").Value = 2250
$ws.Range("K63").Value = 6750
$ws.Range("M63").Value = -6001
$ws.Range("H66").Value = 2250
$ws.Range("I66").Value = 2250
$ws.Range("K66").Value = 20250
$ws.Range("M66").Value = -16506
$ws.Range("H121").Value = 457.125
$ws.Range("J121").Value = 549.6667
$ws.Range("L121").Value = 1649.0001
$ws.Range("N121").Value = -4269.0001
$ws.Range("H122").Value = 2515.8
$ws.Range("I122").Value = 1859.6666
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 16736.9994
$ws.Range("L122").Value = 31500
$ws.Range("M122").Value = -14286.9994
$ws.Range("N122").Value = -36400
$ws.Range("H135").Value = 490
$ws.Range("I135").Value = 490
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4410
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1875
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1262.625
$ws.Range("J97").Value = 1232.2858
$ws.Range("L97").Value = 1232.2858
$ws.Range("N97").Value = -2224.2858
$ws.Range("H102").Value = 4639.8667
$ws.Range("I102").Value = 3780.1
$ws.Range("K102").Value = 3780.1
$ws.Range("M102").Value = -2158.1
$ws.Range("H126").Value = 6134.409
$ws.Range("I126").Value = 5556.643
$ws.Range("J126").Value = 7145.5
$ws.Range("K126").Value = 16669.929
$ws.Range("L126").Value = 21436.5
$ws.Range("M126").Value = -14199.929
$ws.Range("N126").Value = -26376.5
$ws.Range("H132").Value = 127491.25
$ws.Range("I132").Value = 169030
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 507090
$ws.Range("L132").Value = 8625
$ws.Range("M132").Value = -504560
$ws.Range("N132").Value = -13685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2185.3076
$ws.Range("I93").Value = 2246.3635
$ws.Range("J93").Value = 1849.5
$ws.Range("K93").Value = 2246.3635
$ws.Range("L93").Value = 1849.5
$ws.Range("M93").Value = -998.3634999999999
$ws.Range("N93").Value = -4345.5
$ws.Range("H100").Value = 2684.5557
$ws.Range("I100").Value = 2285.5833
$ws.Range("J100").Value = 3482.5
$ws.Range("K100").Value = 2285.5833
$ws.Range("L100").Value = 3482.5
$ws.Range("M100").Value = -1744.5833
$ws.Range("N100").Value = -4564.5
$ws.Range("H127").Value = 101715
$ws.Range("J127").Value = 101715
$ws.Range("L127").Value = 101715
$ws.Range("N127").Value = -111635
$ws.Range("H132").Value = 38319.207
$ws.Range("I132").Value = 40979.16
$ws.Range("K132").Value = 122937.48
$ws.Range("M132").Value = -120407.48
$ws.Range("H133").Value = 88340.5
$ws.Range("J133").Value = 88340.5
$ws.Range("L133").Value = 88340.5
$ws.Range("N133").Value = -93400.5
$ws.Range("H136").Value = 5849.6
$ws.Range("I136").Value = 4187.125
$ws.Range("K136").Value = 12561.375
$ws.Range("M136").Value = -10011.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 761499.5
$ws.Range("J29").Value = 15332.667
$ws.Range("L29").Value = 15332.667
$ws.Range("N29").Value = -15912.667
$ws.Range("H41").Value = 14528
$ws.Range("J41").Value = 14528
$ws.Range("L41").Value = 14528
$ws.Range("N41").Value = -15308
$ws.Range("H132").Value = 113653.445
$ws.Range("I132").Value = 113653.445
$ws.Range("K132").Value = 340960.335
$ws.Range("M132").Value = -338430.335
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
$ws.Range("H136").Value = 5287.3486
$ws.Range("I136").Value = 6065.5483
$ws.Range("K136").Value = 18196.6449
$ws.Range("M136").Value = -15646.6449

